# Rename the asset labels in column A (rows 2-8) from Bloomberg tickers
# to plain-English asset-class names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Cash"
$ws.Range("A3").Value = "EU Flot"
$ws.Range("A4").Value = "EU Equity"
$ws.Range("A5").Value = "US Equity"
$ws.Range("A6").Value = "Greek Gov"
$ws.Range("A7").Value = "EU Corps"
$ws.Range("A8").Value = "EU Gov"

# Update the "Opt Portfolio with View" column (D) so it matches the
# recalculated "Opt Portfolio" column (C) values.
$ws.Range("D2").Value = 0.1535818517403586
$ws.Range("D3").Value = 0.1522119952635631
$ws.Range("D4").Value = 0.1315750985927389
$ws.Range("D5").Value = 0.1316127758639908
$ws.Range("D6").Value = 0.1315393293010641
$ws.Range("D7").Value = 0.149528721328849
$ws.Range("D8").Value = 0.1499502279094355
